# Updates Leve market-price / profit figures across all 8 crafting-job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect a refreshed market-board
# snapshot pulled by the scheduled Sheets runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 480.1
$ws.Range("J17").Value = 480.1
$ws.Range("L17").Value = 1440.3
$ws.Range("N17").Value = -1776.3
$ws.Range("H33").Value = 1286.2
$ws.Range("I33").Value = 1100.6
$ws.Range("J33").Value = 2028.6
$ws.Range("K33").Value = 1100.6
$ws.Range("L33").Value = 2028.6
$ws.Range("M33").Value = -871.5999999999999
$ws.Range("N33").Value = -2486.6
$ws.Range("H51").Value = 38610.39
$ws.Range("J51").Value = 38610.39
$ws.Range("L51").Value = 38610.39
$ws.Range("N51").Value = -39578.39
$ws.Range("H58").Value = 11768.333
$ws.Range("I58").Value = 1015
$ws.Range("J58").Value = 12536.429
$ws.Range("K58").Value = 3045
$ws.Range("L58").Value = 37609.287
$ws.Range("M58").Value = -2895
$ws.Range("N58").Value = -37909.287
$ws.Range("H132").Value = 2766.5334
$ws.Range("I132").Value = 2082.3928
$ws.Range("K132").Value = 6247.178400000001
$ws.Range("M132").Value = -3717.178400000001
$ws.Range("H137").Value = 9680.25
$ws.Range("I137").Value = 4794
$ws.Range("J137").Value = 11309
$ws.Range("K137").Value = 14382
$ws.Range("L137").Value = 33927
$ws.Range("M137").Value = -11832
$ws.Range("N137").Value = -39027
$ws.Range("H138").Value = 3474.0564
$ws.Range("J138").Value = 3774.6072
$ws.Range("L138").Value = 11323.8216
$ws.Range("N138").Value = -21603.8216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7823909
$ws.Range("I32").Value = 8627127
$ws.Range("K32").Value = 8627127
$ws.Range("M32").Value = -8626840
$ws.Range("H88").Value = 2423.7273
$ws.Range("I88").Value = 3179.6
$ws.Range("K88").Value = 3179.6
$ws.Range("M88").Value = -2773.6
$ws.Range("H91").Value = 2423.7273
$ws.Range("I91").Value = 3179.6
$ws.Range("K91").Value = 3179.6
$ws.Range("M91").Value = -1775.6
$ws.Range("H97").Value = 1850.4117
$ws.Range("I97").Value = 1482.6666
$ws.Range("J97").Value = 2733
$ws.Range("K97").Value = 1482.6666
$ws.Range("L97").Value = 2733
$ws.Range("M97").Value = -986.6666
$ws.Range("N97").Value = -3725
$ws.Range("H101").Value = 97994.5
$ws.Range("J101").Value = 97994.5
$ws.Range("L101").Value = 97994.5
$ws.Range("N101").Value = -104484.5
$ws.Range("H132").Value = 5543.352
$ws.Range("I132").Value = 2340.15
$ws.Range("K132").Value = 7020.450000000001
$ws.Range("M132").Value = -4490.450000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1557.125
$ws.Range("I86").Value = 1594.2667
$ws.Range("K86").Value = 1594.2667
$ws.Range("M86").Value = -471.2666999999999
$ws.Range("H89").Value = 1557.125
$ws.Range("I89").Value = 1594.2667
$ws.Range("K89").Value = 7971.3335
$ws.Range("M89").Value = -2355.3335
$ws.Range("H107").Value = 1610
$ws.Range("I107").Value = 1275
$ws.Range("K107").Value = 1275
$ws.Range("M107").Value = 645
$ws.Range("H126").Value = 59999.5
$ws.Range("J126").Value = 59999.5
$ws.Range("L126").Value = 59999.5
$ws.Range("N126").Value = -69879.5
$ws.Range("H130").Value = 96570.57000000001
$ws.Range("J130").Value = 96570.57000000001
$ws.Range("L130").Value = 96570.57000000001
$ws.Range("N130").Value = -106610.57

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 139.6
$ws.Range("I7").Value = 139.6
$ws.Range("K7").Value = 139.6
$ws.Range("M7").Value = -26.59999999999999
$ws.Range("H16").Value = 18252.5
$ws.Range("I16").Value = 22670
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 22670
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -22383
$ws.Range("N16").Value = -5574
$ws.Range("H31").Value = 549718.2
$ws.Range("I31").Value = 8877.429
$ws.Range("K31").Value = 8877.429
$ws.Range("M31").Value = -8582.429
$ws.Range("H34").Value = 549718.2
$ws.Range("I34").Value = 8877.429
$ws.Range("K34").Value = 8877.429
$ws.Range("M34").Value = -8675.429
$ws.Range("H95").Value = 2700
$ws.Range("J95").Value = 2700
$ws.Range("L95").Value = 2700
$ws.Range("N95").Value = -8192
$ws.Range("H113").Value = 18252.5
$ws.Range("I113").Value = 22670
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 22670
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -20500
$ws.Range("N113").Value = -9340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2173.3572
$ws.Range("I68").Value = 2749.75
$ws.Range("J68").Value = 2077.2917
$ws.Range("K68").Value = 8249.25
$ws.Range("L68").Value = 6231.875100000001
$ws.Range("M68").Value = -7438.25
$ws.Range("N68").Value = -7853.875100000001
$ws.Range("H71").Value = 2173.3572
$ws.Range("I71").Value = 2749.75
$ws.Range("J71").Value = 2077.2917
$ws.Range("K71").Value = 24747.75
$ws.Range("L71").Value = 18695.6253
$ws.Range("M71").Value = -20691.75
$ws.Range("N71").Value = -26807.6253
$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 45000
$ws.Range("N101").Value = -49868

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H118").Value = 35510
$ws.Range("J118").Value = 35510
$ws.Range("L118").Value = 35510
$ws.Range("N118").Value = -38824

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 66667136
$ws.Range("J55").Value = 262.5
$ws.Range("L55").Value = 262.5
$ws.Range("N55").Value = -608.5
$ws.Range("H93").Value = 58825444
$ws.Range("I93").Value = 90910500
$ws.Range("J93").Value = 2841.6667
$ws.Range("K93").Value = 90910500
$ws.Range("L93").Value = 2841.6667
$ws.Range("M93").Value = -90909252
$ws.Range("N93").Value = -5337.6667
$ws.Range("H132").Value = 776650.9399999999
$ws.Range("I132").Value = 1255807.8
$ws.Range("K132").Value = 3767423.4
$ws.Range("M132").Value = -3764893.4
$ws.Range("H136").Value = 318643.28
$ws.Range("J136").Value = 303751
$ws.Range("L136").Value = 911253
$ws.Range("N136").Value = -916353

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 25007000
$ws.Range("I62").Value = 7998
$ws.Range("J62").Value = 28578286
$ws.Range("K62").Value = 7998
$ws.Range("L62").Value = 28578286
$ws.Range("M62").Value = -7374
$ws.Range("N62").Value = -28579534
$ws.Range("H65").Value = 25007000
$ws.Range("I65").Value = 7998
$ws.Range("J65").Value = 28578286
$ws.Range("K65").Value = 39990
$ws.Range("L65").Value = 142891430
$ws.Range("M65").Value = -36870
$ws.Range("N65").Value = -142897670
$ws.Range("H132").Value = 246781.61
$ws.Range("I132").Value = 1872.6
$ws.Range("K132").Value = 5617.799999999999
$ws.Range("M132").Value = -3087.799999999999
